$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Anonymize the EC2 instance id in row 2
$ws.Range("A2").Value = "i-xxxxx"

# Add new row 3 for the RDS resource
$ws.Range("A3").Value = "arn:aws:rds:ap-southeast-1:xxxx:db:simple-db"
$ws.Range("B3").Value = "simple-db"
$ws.Range("C3").Value = "RDS"
$ws.Range("D3").Value = "DB"
$ws.Range("E3").Value = "ap-southeast-1"
$ws.Range("F3").Value = "uat"
$ws.Range("G3").Value = "Moodle"
$ws.Range("H3").Value = "'true"
